$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '64.428.60'
Set-TextValue $ws.Range("E2") '  +0.58%  '

Set-TextValue $ws.Range("D3") '3.141.76'
Set-TextValue $ws.Range("E3") '  +0.01%  '

Set-TextValue $ws.Range("E4") '  -0.09%  '

Set-TextValue $ws.Range("D5") '608.88'
Set-TextValue $ws.Range("E5") '  +0.43%  '

Set-TextValue $ws.Range("D6") '144.02'
Set-TextValue $ws.Range("E6") '  -1.66%  '

Set-TextValue $ws.Range("E7") '  -0.10%  '

Set-TextValue $ws.Range("D8") '3.139.19'
Set-TextValue $ws.Range("E8") '  +0.15%  '

Set-TextValue $ws.Range("D9") '0.526'
Set-TextValue $ws.Range("E9") '  +0.67%  '

Set-TextValue $ws.Range("D10") '0.151'
Set-TextValue $ws.Range("E10") '  +0.40%  '

Set-TextValue $ws.Range("E11") '  -1.81%  '

Set-TextValue $ws.Range("D12") '0.471'
Set-TextValue $ws.Range("E12") '  -0.23%  '

Set-TextValue $ws.Range("D13") '0.0000256'
Set-TextValue $ws.Range("E13") '  +2.80%  '

Set-TextValue $ws.Range("D14") '35.53'
Set-TextValue $ws.Range("E14") '  +0.45%  '

Set-TextValue $ws.Range("D15") '3.651.14'
Set-TextValue $ws.Range("E15") '  -0.14%  '

Set-TextValue $ws.Range("E16") '  +2.48%  '

Set-TextValue $ws.Range("D17") '64.296.57'
Set-TextValue $ws.Range("E17") '  +0.32%  '

Set-TextValue $ws.Range("D18") '3.154.33'
Set-TextValue $ws.Range("E18") '  +0.35%  '

Set-TextValue $ws.Range("D19") '6.88'
Set-TextValue $ws.Range("E19") '  +0.22%  '

Set-TextValue $ws.Range("D20") '477.10'
Set-TextValue $ws.Range("E20") '  -0.12%  '

Set-TextValue $ws.Range("D21") '14.87'
Set-TextValue $ws.Range("E21") '  +0.20%  '

Set-TextValue $ws.Range("D22") '0.721'
Set-TextValue $ws.Range("E22") '  +2.50%  '

Set-TextValue $ws.Range("D23") '7.78'
Set-TextValue $ws.Range("E23") '  +0.11%  '

Set-TextValue $ws.Range("D24") '85.28'
Set-TextValue $ws.Range("E24") '  +2.37%  '

Set-TextValue $ws.Range("E25") '  -0.91%  '

Set-TextValue $ws.Range("E26") '  -0.08%  '

Set-TextValue $ws.Range("D27") '2.78'
Set-TextValue $ws.Range("E27") '  -3.21%  '

Set-TextValue $ws.Range("D28") '8.49'
Set-TextValue $ws.Range("E28") '  +1.37%  '

Set-TextValue $ws.Range("D29") '7.35'
Set-TextValue $ws.Range("E29") '  +9.44%  '

Set-TextValue $ws.Range("E30") '  +3.31%  '

Set-TextValue $ws.Range("D31") '2.06'
Set-TextValue $ws.Range("E31") '  -5.12%  '

Set-TextValue $ws.Range("B32") 'EthereumClassic'
Set-TextValue $ws.Range("C32") 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range("D32") '26.98'
Set-TextValue $ws.Range("E32") '  +3.60%  '

Set-TextValue $ws.Range("B33") 'FirstDigitalUSD'
Set-TextValue $ws.Range("C33") 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws.Range("D33") '1.00'
Set-TextValue $ws.Range("E33") '  -0.10%  '

Set-TextValue $ws.Range("D34") '2.65'
Set-TextValue $ws.Range("E34") '  -3.10%  '

Set-TextValue $ws.Range("E35") '  +1.04%  '

Set-TextValue $ws.Range("D36") '5.98'
Set-TextValue $ws.Range("E36") '  +0.85%  '

Set-TextValue $ws.Range("B37") 'PEPE'
Set-TextValue $ws.Range("C37") 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue $ws.Range("D37") '0.0₃0757'
Set-TextValue $ws.Range("E37") '  +3.32%  '

Set-TextValue $ws.Range("B38") 'OKB'
Set-TextValue $ws.Range("C38") 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range("D38") '52.53'
Set-TextValue $ws.Range("E38") '  -2.73%  '

Set-TextValue $ws.Range("D39") '3.02'
Set-TextValue $ws.Range("E39") '  +3.35%  '

Set-TextValue $ws.Range("D40") '445.61'
Set-TextValue $ws.Range("E40") '  -3.70%  '

Set-TextValue $ws.Range("D41") '0.0394'
Set-TextValue $ws.Range("E41") '  +0.51%  '

Set-TextValue $ws.Range("D42") '0.119'
Set-TextValue $ws.Range("E42") '  +0.86%  '

Set-TextValue $ws.Range("D43") '8.30'
Set-TextValue $ws.Range("E43") '  -0.98%  '

Set-TextValue $ws.Range("D44") '2.895.46'
Set-TextValue $ws.Range("E44") '  +2.05%  '

Set-TextValue $ws.Range("D45") '0.264'
Set-TextValue $ws.Range("E45") '  -0.34%  '

Set-TextValue $ws.Range("D46") '2.25'
Set-TextValue $ws.Range("E46") '  -0.08%  '

Set-TextValue $ws.Range("D47") '2.42'
Set-TextValue $ws.Range("E47") '  +5.47%  '

Set-TextValue $ws.Range("B48") 'InjectiveProtocol'
Set-TextValue $ws.Range("C48") 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range("D48") '26.36'
Set-TextValue $ws.Range("E48") '  -0.06%  '

Set-TextValue $ws.Range("B49") 'USDe'
Set-TextValue $ws.Range("C49") 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue $ws.Range("D49") '0.999'
Set-TextValue $ws.Range("E49") '  +0.03%  '

Set-TextValue $ws.Range("E50") '  -0.43%  '

Set-TextValue $ws.Range("B51") 'Arweave'
Set-TextValue $ws.Range("C51") 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
Set-TextValue $ws.Range("D51") '33.91'
Set-TextValue $ws.Range("E51") '  +7.92%  '
